# BOM.xlsx update: add a new BOM line item (row 6) for the WS2812B 8x8 LED matrix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row: No. = 4, Name = "Mạch WS2812B ma trận 8*8", Quantity = 1
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Mạch WS2812B ma trận 8*8"
$ws.Range("C6").Value = 1

# Leave the selection where the author ended up after entering the row.
$ws.Range("B7").Select() | Out-Null
